$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.657.60"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").Value = "1.983.86"
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'246.18"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("E6").Value = "  -4.45%  "
$ws.Range("D7").Value = "'58.31"
$ws.Range("E7").Value = "  +6.77%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'58.74"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'0.362"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "'0.0738"
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("D13").Value = "'0.958"
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("D14").Value = "'14.59"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "2.271.82"
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D17").Value = "1.983.36"
$ws.Range("E17").Value = "  -3.66%  "
$ws.Range("D18").Value = "'18.53"
$ws.Range("E18").Value = "  +8.14%  "
$ws.Range("D19").Value = "35.563.48"
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("D20").Value = "'71.54"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "0.0₃0850"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("D22").Value = "'5.24"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "'233.04"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'2.59"
$ws.Range("E25").Value = "  +21.37%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  -5.29%  "
$ws.Range("D27").Value = "'165.03"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").Value = "'19.26"
$ws.Range("E29").Value = "  -4.47%  "
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("D31").Value = "'4.90"
$ws.Range("E31").Value = "  -3.76%  "
$ws.Range("E32").Value = "  -7.44%  "
$ws.Range("D33").Value = "'0.0953"
$ws.Range("E33").Value = "  +12.74%  "
$ws.Range("D34").Value = "'0.0597"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'2.45"
$ws.Range("E35").Value = "  +10.94%  "
$ws.Range("E36").Value = "  -2.83%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'1.78"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("D39").Value = "'5.46"
$ws.Range("E39").Value = "  +10.03%  "
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "'1.10"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'93.55"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'7.77"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").Value = "'0.0908"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "'16.16"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "1.375.15"
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("D49").Value = "'2.90"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "'46.90"
$ws.Range("E50").Value = "  +2.36%  "
$ws.Range("D51").Value = "'2.28"
$ws.Range("E51").Value = "  -0.05%  "
